# Applies the scheduled market-data refresh to the Chocobo_Profits crafting-class
# sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR): updated currentAveragePrice* /
# LevePrice*/LeveProfit* columns (H:N) per leve row, as produced by the runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1628.2222
$ws.Range("I40").Value = 1612.5
$ws.Range("J40").Value = 1640.8
$ws.Range("K40").Value = 1612.5
$ws.Range("L40").Value = 1640.8
$ws.Range("M40").Value = -1437.5
$ws.Range("N40").Value = -1990.8
$ws.Range("H98").Value = 6177.273
$ws.Range("I98").Value = 5051.2144
$ws.Range("J98").Value = 7345.037
$ws.Range("K98").Value = 5051.2144
$ws.Range("L98").Value = 7345.037
$ws.Range("M98").Value = -3553.2144
$ws.Range("N98").Value = -10341.037
$ws.Range("H112").Value = 1312.4263
$ws.Range("I112").Value = 566.6667
$ws.Range("J112").Value = 1351
$ws.Range("K112").Value = 1700.0001
$ws.Range("L112").Value = 4053
$ws.Range("M112").Value = -592.0001
$ws.Range("N112").Value = -6269
$ws.Range("I116").Value = 1251737.5
$ws.Range("J116").Value = 7034.7896
$ws.Range("K116").Value = 1251737.5
$ws.Range("L116").Value = 7034.7896
$ws.Range("M116").Value = -1248295.5
$ws.Range("N116").Value = -13918.7896
$ws.Range("H122").Value = 6177.273
$ws.Range("I122").Value = 5051.2144
$ws.Range("J122").Value = 7345.037
$ws.Range("K122").Value = 15153.6432
$ws.Range("L122").Value = 22035.111
$ws.Range("M122").Value = -12703.6432
$ws.Range("N122").Value = -26935.111
$ws.Range("H129").Value = 846.24
$ws.Range("I129").Value = 375.5
$ws.Range("J129").Value = 865.8542
$ws.Range("K129").Value = 1126.5
$ws.Range("L129").Value = 2597.5626
$ws.Range("M129").Value = 3873.5
$ws.Range("N129").Value = -12597.5626
$ws.Range("H138").Value = 2675.465
$ws.Range("I138").Value = 984.3333
$ws.Range("J138").Value = 3032.7466
$ws.Range("K138").Value = 2952.9999
$ws.Range("L138").Value = 9098.239799999999
$ws.Range("M138").Value = 2187.0001
$ws.Range("N138").Value = -19378.2398

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3017.5757
$ws.Range("I32").Value = 2780.013
$ws.Range("J32").Value = 3849.0454
$ws.Range("K32").Value = 2780.013
$ws.Range("L32").Value = 3849.0454
$ws.Range("M32").Value = -2493.013
$ws.Range("N32").Value = -4423.0454
$ws.Range("H61").Value = 1353
$ws.Range("I61").Value = 1267.0834
$ws.Range("K61").Value = 1267.0834
$ws.Range("M61").Value = -1055.0834
$ws.Range("H97").Value = 824.08
$ws.Range("I97").Value = 688.8889
$ws.Range("J97").Value = 1171.7142
$ws.Range("K97").Value = 688.8889
$ws.Range("L97").Value = 1171.7142
$ws.Range("M97").Value = -192.8889
$ws.Range("N97").Value = -2163.7142
$ws.Range("H115").Value = 28088.8
$ws.Range("J115").Value = 28088.8
$ws.Range("L115").Value = 28088.8
$ws.Range("N115").Value = -31222.8
$ws.Range("H132").Value = 2716.8235
$ws.Range("I132").Value = 1347.8334
$ws.Range("J132").Value = 6002.4
$ws.Range("K132").Value = 4043.5002
$ws.Range("L132").Value = 18007.2
$ws.Range("M132").Value = -1513.5002
$ws.Range("N132").Value = -23067.2
$ws.Range("H136").Value = 1353
$ws.Range("I136").Value = 1267.0834
$ws.Range("K136").Value = 3801.2502
$ws.Range("M136").Value = -1251.2502
$ws.Range("H139").Value = 41498.93
$ws.Range("J139").Value = 41498.93
$ws.Range("L139").Value = 41498.93
$ws.Range("N139").Value = -51778.93

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 360.22223
$ws.Range("I94").Value = 345.3
$ws.Range("J94").Value = 402.85715
$ws.Range("K94").Value = 345.3
$ws.Range("L94").Value = 402.85715
$ws.Range("M94").Value = 105.7
$ws.Range("N94").Value = -1304.85715
$ws.Range("H134").Value = 3090.75
$ws.Range("I134").Value = 1487.5625
$ws.Range("J134").Value = 9503.5
$ws.Range("K134").Value = 4462.6875
$ws.Range("L134").Value = 28510.5
$ws.Range("M134").Value = -1927.6875
$ws.Range("N134").Value = -33580.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 14715715
$ws.Range("I6").Value = 14715715
$ws.Range("K6").Value = 14715715
$ws.Range("M6").Value = -14715602
$ws.Range("H31").Value = 7524.515
$ws.Range("I31").Value = 2290.2
$ws.Range("K31").Value = 2290.2
$ws.Range("M31").Value = -1995.2
$ws.Range("H34").Value = 7524.515
$ws.Range("I34").Value = 2290.2
$ws.Range("K34").Value = 2290.2
$ws.Range("M34").Value = -2088.2
$ws.Range("H99").Value = 20005560
$ws.Range("I99").Value = 66668268
$ws.Range("J99").Value = 7257.143
$ws.Range("K99").Value = 66668268
$ws.Range("L99").Value = 7257.143
$ws.Range("M99").Value = -66666770
$ws.Range("N99").Value = -10253.143
$ws.Range("H122").Value = 3741.5715
$ws.Range("I122").Value = 1797.75
$ws.Range("J122").Value = 6333.3335
$ws.Range("K122").Value = 5393.25
$ws.Range("L122").Value = 19000.0005
$ws.Range("N122").Value = -23900.0005
$ws.Range("M122").Value = -2943.25
$ws.Range("H126").Value = 20005560
$ws.Range("I126").Value = 66668268
$ws.Range("J126").Value = 7257.143
$ws.Range("K126").Value = 200004804
$ws.Range("L126").Value = 21771.429
$ws.Range("M126").Value = -200002334
$ws.Range("N126").Value = -26711.429
$ws.Range("H132").Value = 2834.1853
$ws.Range("I132").Value = 1395.9445
$ws.Range("J132").Value = 5710.6665
$ws.Range("K132").Value = 4187.833500000001
$ws.Range("L132").Value = 17131.9995
$ws.Range("M132").Value = -1657.833500000001
$ws.Range("N132").Value = -22191.9995
$ws.Range("H134").Value = 5523.207
$ws.Range("I134").Value = 7796.9375
$ws.Range("J134").Value = 2724.7693
$ws.Range("K134").Value = 23390.8125
$ws.Range("L134").Value = 8174.3079
$ws.Range("M134").Value = -20855.8125
$ws.Range("N134").Value = -13244.3079
$ws.Range("H137").Value = 30456.666
$ws.Range("J137").Value = 30456.666
$ws.Range("L137").Value = 30456.666
$ws.Range("N137").Value = -40656.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 299.08334
$ws.Range("I92").Value = 299.08334
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 897.2500200000001
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 350.7499799999999
$ws.Range("N92").ClearContents()
$ws.Range("H132").Value = 2032.4062
$ws.Range("I132").Value = 908.0909
$ws.Range("J132").Value = 2621.3333
$ws.Range("K132").Value = 8172.8181
$ws.Range("L132").Value = 23591.9997
$ws.Range("M132").Value = -5642.8181
$ws.Range("N132").Value = -28651.9997
$ws.Range("H140").Value = 15700
$ws.Range("I140").Value = 30000
$ws.Range("J140").Value = 1400
$ws.Range("K140").Value = 90000
$ws.Range("L140").Value = 4200
$ws.Range("M140").Value = -84820
$ws.Range("N140").Value = -14560
$ws.Range("H141").Value = 8656.429
$ws.Range("I141").Value = 9132.5
$ws.Range("K141").Value = 27397.5
$ws.Range("M141").Value = -22217.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4754.846
$ws.Range("I122").Value = 3781.3
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 11343.9
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -8893.900000000001
$ws.Range("N122").Value = -28900
$ws.Range("H132").Value = 3072.5454
$ws.Range("I132").Value = 1790.3334
$ws.Range("J132").Value = 4611.2
$ws.Range("K132").Value = 5371.0002
$ws.Range("L132").Value = 13833.6
$ws.Range("M132").Value = -2841.0002
$ws.Range("N132").Value = -18893.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1718.4546
$ws.Range("J61").Value = 1716.5
$ws.Range("L61").Value = 1716.5
$ws.Range("N61").Value = -2120.5
$ws.Range("H96").Value = 38800
$ws.Range("J96").Value = 38800
$ws.Range("L96").Value = 38800
$ws.Range("N96").Value = -44292
$ws.Range("H113").Value = 1718.4546
$ws.Range("J113").Value = 1716.5
$ws.Range("L113").Value = 1716.5
$ws.Range("N113").Value = -6056.5
